# Implements the "copySheet" update:
#  - the previously-copied sheet ("Sheet0") is renamed to "copied_sheet" and
#    moved in front of "Sheet1", becoming the first / active tab
#  - a new cell (G13) is written on that sheet containing a single string
#    that mixes several different text formats in one cell
#  - the pre-existing text cells that went through the copy routine now
#    carry their formatting as explicit (per-run) rich text instead of a
#    plain shared string, even though the rendered look is unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet0")

# --- rename + reposition the copied sheet, make it the active tab -------
$ws.Name = "copied_sheet"
$ws.Move($wb.Worksheets.Item(1))

# re-fetch: Move() leaves the old $ws reference pointing at whatever sheet
# now occupies the slot the moved sheet used to be in, so grab a fresh
# handle by name before doing any further work on it
$ws = $wb.Worksheets.Item("copied_sheet")

# --- helper: force explicit per-character rich-text runs -----------------
function Set-RunFont($range, $start, $length, $size, $colorIndex) {
    $chars = $range.Characters($start, $length)
    $chars.Font.Name = "宋体"
    $chars.Font.Size = $size
    $chars.Font.ColorIndex = $colorIndex
}

# --- new cell: several different formats inside a single cell -----------
$g13 = $ws.Range("G13")
$g13.Value = "多种文本在一个Cell里"

Set-RunFont $g13 1 1 12 1
Set-RunFont $g13 2 1 12 1
Set-RunFont $g13 3 1 11 46
Set-RunFont $g13 4 1 11 46
Set-RunFont $g13 5 1 12 46
Set-RunFont $g13 6 1 22 3
Set-RunFont $g13 7 1 12 1
Set-RunFont $g13 12 1 12 1

# --- re-assert the (unchanged) formatting of the pre-existing text cells
# as explicit rich text, matching what the updated copySheet routine
# produces when it re-writes cell styles during the copy -----------------
$a1 = $ws.Range("A1")
Set-RunFont $a1 1 1 12 1
Set-RunFont $a1 2 1 12 1
Set-RunFont $a1 3 1 12 1

$a35 = $ws.Range("A35")
Set-RunFont $a35 1 1 12 1
Set-RunFont $a35 2 1 12 1

$a36 = $ws.Range("A36")
Set-RunFont $a36 1 1 12 1
Set-RunFont $a36 2 1 12 1
Set-RunFont $a36 3 1 12 1

Write-Output "copySheet update applied"
